{"js": "// Update the lattice-multiplication exercise table: each cell holds a\n// multi-line problem (\"AB x CD\", the second factor split apart, a\n// constant \"  ----\" divider, and the two digits of the first factor\n// down the left side), with lines separated by manual line breaks\n// (<w:br/>). The new set of 15 problems replaces the old ones\n// cell-for-cell while keeping the 5x3 grid and per-run formatting\n// (sz 32, and the xml:space=\"preserve\" that the space-padded lines\n// need) untouched - done by replacing each cell body with an OOXML\n// fragment that mirrors the original run/paragraph shape exactly.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// row, col (0-based), equation line, second-factor line, left digit 1,\n// left digit 2 - mirrors the cell-by-cell content dictated by the diff.\nconst cells = [\n  { row: 0, col: 0, eq: \"69 x 93\", second: \"  9    3\", l1: \"6\", l2: \"9\" },\n  { row: 0, col: 1, eq: \"33 x 24\", second: \"  2    4\", l1: \"3\", l2: \"3\" },\n  { row: 0, col: 2, eq: \"83 x 41\", second: \"  4    1\", l1: \"8\", l2: \"3\" },\n  { row: 1, col: 0, eq: \"15 x 93\", second: \"  9    3\", l1: \"1\", l2: \"5\" },\n  { row: 1, col: 1, eq: \"66 x 51\", second: \"  5    1\", l1: \"6\", l2: \"6\" },\n  { row: 1, col: 2, eq: \"31 x 39\", second: \"  3    9\", l1: \"3\", l2: \"1\" },\n  { row: 2, col: 0, eq: \"87 x 91\", second: \"  9    1\", l1: \"8\", l2: \"7\" },\n  { row: 2, col: 1, eq: \"30 x 61\", second: \"  6    1\", l1: \"3\", l2: \"0\" },\n  { row: 2, col: 2, eq: \"18 x 84\", second: \"  8    4\", l1: \"1\", l2: \"8\" },\n  { row: 3, col: 0, eq: \"94 x 45\", second: \"  4    5\", l1: \"9\", l2: \"4\" },\n  { row: 3, col: 1, eq: \"67 x 57\", second: \"  5    7\", l1: \"6\", l2: \"7\" },\n  { row: 3, col: 2, eq: \"69 x 87\", second: \"  8    7\", l1: \"6\", l2: \"9\" },\n  { row: 4, col: 0, eq: \"63 x 29\", second: \"  2    9\", l1: \"6\", l2: \"3\" },\n  { row: 4, col: 1, eq: \"36 x 68\", second: \"  6    8\", l1: \"3\", l2: \"6\" },\n  { row: 4, col: 2, eq: \"38 x 45\", second: \"  4    5\", l1: \"3\", l2: \"8\" },\n];\n\nfunction cellOoxml(c) {\n  return (\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    \"<pkg:xmlData>\" +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    \"<w:body><w:p><w:r><w:rPr><w:sz w:val=\\\"32\\\"/></w:rPr>\" +\n    \"<w:t>\" + c.eq + \"</w:t><w:br/>\" +\n    '<w:t xml:space=\"preserve\">' + c.second + \"</w:t><w:br/>\" +\n    '<w:t xml:space=\"preserve\">  ----</w:t><w:br/>' +\n    \"<w:t>\" + c.l1 + \"|    |</w:t><w:br/>\" +\n    \"<w:t>\" + c.l2 + \"|    |</w:t>\" +\n    \"</w:r></w:p></w:body></w:document>\" +\n    \"</pkg:xmlData></pkg:part></pkg:package>\"\n  );\n}\n\nfor (const c of cells) {\n  const cell = table.getCell(c.row, c.col);\n  cell.body.insertOoxml(cellOoxml(c), Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Update the lattice-multiplication exercise table: each cell holds a\n# multi-line problem (\"AB x CD\", the second factor split apart, a\n# constant \"  ----\" divider, and the two digits of the first factor\n# down the left side), with lines separated by manual line breaks\n# (<w:br/>). The new set of 15 problems replaces the old ones\n# cell-for-cell while keeping the 5x3 grid and per-run formatting\n# (sz 32) untouched.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# Row, Col (1-based), equation line, second-factor line, left digit 1,\n# left digit 2 - mirrors the cell-by-cell content dictated by the diff.\n$cells = @(\n    @{Row=1; Col=1; Eq=\"69 x 93\"; Second=\"  9    3\"; L1=\"6\"; L2=\"9\"}\n    @{Row=1; Col=2; Eq=\"33 x 24\"; Second=\"  2    4\"; L1=\"3\"; L2=\"3\"}\n    @{Row=1; Col=3; Eq=\"83 x 41\"; Second=\"  4    1\"; L1=\"8\"; L2=\"3\"}\n    @{Row=2; Col=1; Eq=\"15 x 93\"; Second=\"  9    3\"; L1=\"1\"; L2=\"5\"}\n    @{Row=2; Col=2; Eq=\"66 x 51\"; Second=\"  5    1\"; L1=\"6\"; L2=\"6\"}\n    @{Row=2; Col=3; Eq=\"31 x 39\"; Second=\"  3    9\"; L1=\"3\"; L2=\"1\"}\n    @{Row=3; Col=1; Eq=\"87 x 91\"; Second=\"  9    1\"; L1=\"8\"; L2=\"7\"}\n    @{Row=3; Col=2; Eq=\"30 x 61\"; Second=\"  6    1\"; L1=\"3\"; L2=\"0\"}\n    @{Row=3; Col=3; Eq=\"18 x 84\"; Second=\"  8    4\"; L1=\"1\"; L2=\"8\"}\n    @{Row=4; Col=1; Eq=\"94 x 45\"; Second=\"  4    5\"; L1=\"9\"; L2=\"4\"}\n    @{Row=4; Col=2; Eq=\"67 x 57\"; Second=\"  5    7\"; L1=\"6\"; L2=\"7\"}\n    @{Row=4; Col=3; Eq=\"69 x 87\"; Second=\"  8    7\"; L1=\"6\"; L2=\"9\"}\n    @{Row=5; Col=1; Eq=\"63 x 29\"; Second=\"  2    9\"; L1=\"6\"; L2=\"3\"}\n    @{Row=5; Col=2; Eq=\"36 x 68\"; Second=\"  6    8\"; L1=\"3\"; L2=\"6\"}\n    @{Row=5; Col=3; Eq=\"38 x 45\"; Second=\"  4    5\"; L1=\"3\"; L2=\"8\"}\n)\n\nforeach ($c in $cells) {\n    $cell = $t.Cell($c.Row, $c.Col)\n\n    $xmlFragment = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n        '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n        '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n        '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n        '<w:body><w:p><w:r><w:rPr><w:sz w:val=\"32\"/></w:rPr>' +\n        \"<w:t>$($c.Eq)</w:t><w:br/>\" +\n        \"<w:t xml:space=`\"preserve`\">$($c.Second)</w:t><w:br/>\" +\n        '<w:t xml:space=\"preserve\">  ----</w:t><w:br/>' +\n        \"<w:t>$($c.L1)|    |</w:t><w:br/>\" +\n        \"<w:t>$($c.L2)|    |</w:t>\" +\n        '</w:r></w:p></w:body></w:document>' +\n        '</pkg:xmlData></pkg:part></pkg:package>'\n\n    $cell.Range.InsertXML($xmlFragment)\n}\n\nWrite-Output \"Updated $($cells.Count) lattice-multiplication cells\"\n"}
